$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1152.25
$ws.Range("I2").Value = 109.17647
$ws.Range("J2").Value = 3685.4285
$ws.Range("K2").Value = 109.17647
$ws.Range("L2").Value = 3685.4285
$ws.Range("M2").Value = 3.823530000000005
$ws.Range("N2").Value = -3911.4285
$ws.Range("H32").Value = 924.4167
$ws.Range("I32").Value = 824.75
$ws.Range("K32").Value = 824.75
$ws.Range("M32").Value = -498.75
$ws.Range("H43").Value = 39890
$ws.Range("I43").Value = 1850
$ws.Range("J43").Value = 86383.336
$ws.Range("K43").Value = 1850
$ws.Range("L43").Value = 86383.336
$ws.Range("M43").Value = -1781
$ws.Range("N43").Value = -86521.336
$ws.Range("H92").Value = 294.55554
$ws.Range("I92").Value = 325.83334
$ws.Range("K92").Value = 325.83334
$ws.Range("M92").Value = 922.16666
$ws.Range("H98").Value = 1725.1111
$ws.Range("I98").Value = 1500.0834
$ws.Range("K98").Value = 1500.0834
$ws.Range("M98").Value = -2.083399999999983
$ws.Range("H100").Value = 5797.8
$ws.Range("I100").Value = 4998.3335
$ws.Range("J100").Value = 6997
$ws.Range("K100").Value = 4998.3335
$ws.Range("L100").Value = 6997
$ws.Range("M100").Value = -4457.3335
$ws.Range("N100").Value = -8079
$ws.Range("H106").Value = 3766.1667
$ws.Range("I106").Value = 3766.1667
$ws.Range("K106").Value = 3766.1667
$ws.Range("M106").Value = -3135.1667
$ws.Range("H113").Value = 4143.857
$ws.Range("J113").Value = 4566.25
$ws.Range("L113").Value = 4566.25
$ws.Range("N113").Value = -11074.25
$ws.Range("H122").Value = 1725.1111
$ws.Range("I122").Value = 1500.0834
$ws.Range("K122").Value = 4500.2502
$ws.Range("M122").Value = -2050.2502
$ws.Range("H132").Value = 5596.4287
$ws.Range("I132").Value = 6495.8335
$ws.Range("J132").Value = 200
$ws.Range("K132").Value = 19487.5005
$ws.Range("L132").Value = 600
$ws.Range("M132").Value = -16957.5005
$ws.Range("N132").Value = -5660
$ws.Range("H134").Value = 78000
$ws.Range("J134").Value = 78000
$ws.Range("L134").Value = 78000
$ws.Range("N134").Value = -88140
$ws.Range("H137").Value = 1532.9166
$ws.Range("I137").Value = 1532.9166
$ws.Range("K137").Value = 4598.7498
$ws.Range("M137").Value = -2048.7498
$ws.Range("H138").Value = 1787.375
$ws.Range("I138").Value = 724.25
$ws.Range("K138").Value = 2172.75
$ws.Range("M138").Value = 2967.25

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 21000
$ws.Range("J37").Value = 23750
$ws.Range("L37").Value = 23750
$ws.Range("N37").Value = -24296
$ws.Range("H45").Value = 1994.7
$ws.Range("I45").Value = 1994.7
$ws.Range("K45").Value = 1994.7
$ws.Range("M45").Value = -1617.7
$ws.Range("H110").Value = 2007.1111
$ws.Range("I110").Value = 510.83334
$ws.Range("K110").Value = 510.83334
$ws.Range("M110").Value = 1534.16666
$ws.Range("H132").Value = 1477.8
$ws.Range("I132").Value = 1519.8889
$ws.Range("K132").Value = 4559.6667
$ws.Range("M132").Value = -2029.6667

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2516.8
$ws.Range("I86").Value = 1814.8334
$ws.Range("K86").Value = 1814.8334
$ws.Range("M86").Value = -691.8334
$ws.Range("H89").Value = 2516.8
$ws.Range("I89").Value = 1814.8334
$ws.Range("K89").Value = 9074.166999999999
$ws.Range("M89").Value = -3458.166999999999
$ws.Range("H94").Value = 2921.111
$ws.Range("I94").Value = 2921.111
$ws.Range("K94").Value = 2921.111
$ws.Range("M94").Value = -2470.111
$ws.Range("H99").Value = 2999
$ws.Range("I99").Value = 2999
$ws.Range("J99").Value = 2999
$ws.Range("K99").Value = 2999
$ws.Range("L99").Value = 2999
$ws.Range("M99").Value = -1501
$ws.Range("N99").Value = -5995
$ws.Range("H105").Value = 1015.5
$ws.Range("I105").Value = 787.6667
$ws.Range("J105").Value = 1699
$ws.Range("K105").Value = 787.6667
$ws.Range("L105").Value = 1699
$ws.Range("M105").Value = 959.3333
$ws.Range("N105").Value = -5193
$ws.Range("H134").Value = 9192.8125
$ws.Range("I134").Value = 8826.817999999999
$ws.Range("K134").Value = 26480.454
$ws.Range("M134").Value = -23945.454

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 36666.668
$ws.Range("J39").Value = 52000
$ws.Range("L39").Value = 52000
$ws.Range("N39").Value = -52782
$ws.Range("H49").Value = 36666.668
$ws.Range("J49").Value = 52000
$ws.Range("L49").Value = 52000
$ws.Range("N49").Value = -52364
$ws.Range("H58").Value = 2996
$ws.Range("I58").Value = 2996
$ws.Range("K58").Value = 2996
$ws.Range("M58").Value = -2793
$ws.Range("H105").Value = 663
$ws.Range("I105").Value = 663
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 663
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 1084
$ws.Range("H136").Value = 2996
$ws.Range("I136").Value = 2996
$ws.Range("K136").Value = 8988
$ws.Range("M136").Value = -6438
$ws.Range("N105").Delete()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 33.8
$ws.Range("I38").Value = 30
$ws.Range("K38").Value = 90
$ws.Range("M38").Value = 257
$ws.Range("H116").Value = 3165.3333
$ws.Range("I116").Value = 2750
$ws.Range("K116").Value = 8250
$ws.Range("M116").Value = -4808

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 319.52942
$ws.Range("I2").Value = 361.46155
$ws.Range("K2").Value = 361.46155
$ws.Range("M2").Value = -248.46155
$ws.Range("H80").Value = 2602.0715
$ws.Range("I80").Value = 1375
$ws.Range("J80").Value = 3522.375
$ws.Range("K80").Value = 1375
$ws.Range("L80").Value = 3522.375
$ws.Range("M80").Value = -377
$ws.Range("N80").Value = -5518.375
$ws.Range("H83").Value = 2602.0715
$ws.Range("I83").Value = 1375
$ws.Range("J83").Value = 3522.375
$ws.Range("K83").Value = 6875
$ws.Range("L83").Value = 17611.875
$ws.Range("M83").Value = -1883
$ws.Range("N83").Value = -27595.875
$ws.Range("H113").Value = 1579.6
$ws.Range("I113").Value = 1579.6
$ws.Range("K113").Value = 1579.6
$ws.Range("M113").Value = 590.4000000000001
$ws.Range("H132").Value = 4875
$ws.Range("I132").Value = 4875
$ws.Range("K132").Value = 14625
$ws.Range("M132").Value = -12095

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2554.3635
$ws.Range("I7").Value = 2137.25
$ws.Range("J7").Value = 3666.6667
$ws.Range("K7").Value = 2137.25
$ws.Range("L7").Value = 3666.6667
$ws.Range("M7").Value = -2025.25
$ws.Range("N7").Value = -3890.6667
$ws.Range("H22").Value = 1218.5
$ws.Range("I22").Value = 990
$ws.Range("K22").Value = 990
$ws.Range("M22").Value = -695
$ws.Range("H27").Value = 1218.5
$ws.Range("I27").Value = 990
$ws.Range("K27").Value = 990
$ws.Range("M27").Value = -883
$ws.Range("H46").Value = 2799.3333
$ws.Range("I46").Value = 2181
$ws.Range("K46").Value = 2181
$ws.Range("M46").Value = -1993
$ws.Range("H93").Value = 1002
$ws.Range("J93").Value = 1002
$ws.Range("L93").Value = 1002
$ws.Range("N93").Value = -3498
$ws.Range("H126").Value = 2554.3635
$ws.Range("I126").Value = 2137.25
$ws.Range("J126").Value = 3666.6667
$ws.Range("K126").Value = 6411.75
$ws.Range("L126").Value = 11000.0001
$ws.Range("M126").Value = -3941.75
$ws.Range("N126").Value = -15940.0001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("H92").Value = 8333.333000000001
$ws.Range("J92").Value = 8333.333000000001
$ws.Range("L92").Value = 8333.333000000001
$ws.Range("N92").Value = -13325.333
$ws.Range("H122").Value = 778.0714
$ws.Range("I122").Value = 778.0714
$ws.Range("K122").Value = 2334.2142
$ws.Range("M122").Value = 115.7857999999997
$ws.Range("H126").Value = 2027.1765
$ws.Range("I126").Value = 1747.4286
$ws.Range("J126").Value = 3332.6667
$ws.Range("K126").Value = 5242.2858
$ws.Range("L126").Value = 9998.000100000001
$ws.Range("M126").Value = -2772.2858
$ws.Range("N126").Value = -14938.0001
$ws.Range("M62").Delete()
$ws.Range("M65").Delete()
